$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Merge the run sequence "," + " visa " + "também " +
# "otimizar o espaço nas academias usando " into a single run reading
# ", visa também otimizar o espaço nas academias usando ", while leaving the
# neighbouring runs ("sensores", ", para acomodar ... usuários" and the
# closing ".") untouched.
#
# The editing surface only exposes text-level operations (Find/Replace,
# Range.Text, formatting, ...). Any text replacement inside a paragraph
# causes the engine to recompute run boundaries by merging maximal runs of
# identical character formatting (w:rPr) around the edited text. Since
# "sensores" (and the remaining tail of the sentence) share that exact same
# formatting, a naive replace would also swallow them into the merged run.
# To prevent that we briefly give the following runs a distinguishing
# (and otherwise invisible) formatting nudge, perform the text merge, and
# then revert the nudge - toggling character formatting alone does not
# trigger the run-merging pass, so the boundaries end up exactly where we
# want them.
# ---------------------------------------------------------------------------

# 1) Anchor on the unique bold phrase right before the text we need to touch.
$anchor = $d.Content
[void]$anchor.Find.Execute("superlotação e falta de otimização de tempo nos treinos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterAnchor = $anchor.End

# 2) Locate the four-run phrase that must become a single run.
$mergeRange = $d.Range($afterAnchor, $d.Content.End)
[void]$mergeRange.Find.Execute(", visa também otimizar o espaço nas academias usando ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $mergeRange.Start
$mergeEnd = $mergeRange.End

# 3) Locate "sensores", immediately following the phrase above.
$sensoresRange = $d.Range($mergeEnd, $d.Content.End)
[void]$sensoresRange.Find.Execute("sensores", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sensoresStart = $sensoresRange.Start
$sensoresEnd = $sensoresRange.End

# 4) Locate the closing period that ends the sentence / paragraph.
$dotRange = $d.Range($sensoresEnd, $d.Content.End)
[void]$dotRange.Find.Execute("bem-estar dos usuários.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotEnd = $dotRange.End
$dotStart = $dotEnd - 1

# Protect the boundaries right after the merge target so the upcoming
# replace does not cascade past them.
$protectSensores = $d.Range($sensoresStart, $sensoresEnd)
$protectSensores.Bold = 1

$protectDot = $d.Range($dotStart, $dotEnd)
$protectDot.Bold = 1

# Perform the actual merge: replace the four-run phrase with the identical
# text so the engine rebuilds it as a single run.
$editRange = $d.Range($mergeStart, $mergeEnd)
[void]$editRange.Find.Execute(", visa também otimizar o espaço nas academias usando ", $true, $false, $false, $false, $false, $true, 1, $false, ", visa também otimizar o espaço nas academias usando ", 2)

# Undo the temporary protection formatting; toggling formatting alone does
# not re-trigger run merging, so this restores the original look without
# disturbing the run boundaries we just fixed in place.
$protectSensores2 = $d.Range($sensoresStart, $sensoresEnd)
$protectSensores2.Bold = 0

$protectDot2 = $d.Range($dotStart, $dotEnd)
$protectDot2.Bold = 0

$d.Save()
